$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8546833333333335
$ws.Range("H2").Value = 2.56405
$ws.Range("I2").Value = 0.3097546281380014
$ws.Range("J2").Value = 0.3097546281380015
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.021894333333333
$ws.Range("N2").Value = 9.065683
$ws.Range("O2").Value = 0.1464771679819186
$ws.Range("P2").Value = 0.1464771679819185
$ws.Range("Q2").Value = 2.582762721794445
$ws.Range("R2").Value = 23.24486449615
$ws.Range("S2").Value = 0.04537198069894675
$ws.Range("T2").Value = 0.04537198069894675
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8546833333333335
$ws.Range("H3").Value = 2.56405
$ws.Range("I3").Value = 0.3097546281380014
$ws.Range("J3").Value = 0.3097546281380015
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.799695333333333
$ws.Range("N3").Value = 17.399086
$ws.Range("O3").Value = 0.2811226515149324
$ws.Range("P3").Value = 0.2811226515149324
$ws.Range("Q3").Value = 4.956902939811112
$ws.Range("R3").Value = 44.61212645830001
$ws.Range("S3").Value = 0.08707904238117686
$ws.Range("T3").Value = 0.08707904238117688
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8546833333333335
$ws.Range("H4").Value = 2.56405
$ws.Range("I4").Value = 0.3097546281380014
$ws.Range("J4").Value = 0.3097546281380015
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.259039333333333
$ws.Range("N4").Value = 18.777118
$ws.Range("O4").Value = 0.303387959572633
$ws.Range("P4").Value = 0.303387959572633
$ws.Range("Q4").Value = 5.349496600877778
$ws.Range("R4").Value = 48.1454694079
$ws.Range("S4").Value = 0.09397582459896794
$ws.Range("T4").Value = 0.09397582459896796
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8546833333333335
$ws.Range("H5").Value = 2.56405
$ws.Range("I5").Value = 0.3097546281380014
$ws.Range("J5").Value = 0.3097546281380015
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.549851333333334
$ws.Range("N5").Value = 16.649554
$ws.Range("O5").Value = 0.2690122209305161
$ws.Range("P5").Value = 0.2690122209305161
$ws.Range("Q5").Value = 4.743365437077779
$ws.Range("R5").Value = 42.69028893370001
$ws.Range("S5").Value = 0.0833277804589099
$ws.Range("T5").Value = 0.0833277804589099
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.251276
$ws.Range("H6").Value = 3.753828
$ws.Range("I6").Value = 0.453487879032787
$ws.Range("J6").Value = 0.4534878790327871
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.021894333333333
$ws.Range("N6").Value = 9.065683
$ws.Range("O6").Value = 0.1464771679819186
$ws.Range("P6").Value = 0.1464771679819185
$ws.Range("Q6").Value = 3.781223853836
$ws.Range("R6").Value = 34.031014684524
$ws.Range("S6").Value = 0.06642562023484951
$ws.Range("T6").Value = 0.06642562023484951
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.251276
$ws.Range("H7").Value = 3.753828
$ws.Range("I7").Value = 0.453487879032787
$ws.Range("J7").Value = 0.4534878790327871
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.799695333333333
$ws.Range("N7").Value = 17.399086
$ws.Range("O7").Value = 0.2811226515149324
$ws.Range("P7").Value = 0.2811226515149324
$ws.Range("Q7").Value = 7.257019577912001
$ws.Range("R7").Value = 65.31317620120801
$ws.Range("S7").Value = 0.12748571498358
$ws.Range("T7").Value = 0.1274857149835801
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.251276
$ws.Range("H8").Value = 3.753828
$ws.Range("I8").Value = 0.453487879032787
$ws.Range("J8").Value = 0.4534878790327871
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.259039333333333
$ws.Range("N8").Value = 18.777118
$ws.Range("O8").Value = 0.303387959572633
$ws.Range("P8").Value = 0.303387959572633
$ws.Range("Q8").Value = 7.831785700856
$ws.Range("R8").Value = 70.48607130770399
$ws.Range("S8").Value = 0.1375827623106783
$ws.Range("T8").Value = 0.1375827623106783
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.251276
$ws.Range("H9").Value = 3.753828
$ws.Range("I9").Value = 0.453487879032787
$ws.Range("J9").Value = 0.4534878790327871
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.549851333333334
$ws.Range("N9").Value = 16.649554
$ws.Range("O9").Value = 0.2690122209305161
$ws.Range("P9").Value = 0.2690122209305161
$ws.Range("Q9").Value = 6.944395776968001
$ws.Range("R9").Value = 62.49956199271201
$ws.Range("S9").Value = 0.1219937815036793
$ws.Range("T9").Value = 0.1219937815036793
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1689346666666667
$ws.Range("H10").Value = 0.506804
$ws.Range("I10").Value = 0.06122536009783416
$ws.Range("J10").Value = 0.06122536009783416
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.021894333333333
$ws.Range("N10").Value = 9.065683
$ws.Range("O10").Value = 0.1464771679819186
$ws.Range("P10").Value = 0.1464771679819185
$ws.Range("Q10").Value = 0.5105027119035556
$ws.Range("R10").Value = 4.594524407132
$ws.Range("S10").Value = 0.008968117355803908
$ws.Range("T10").Value = 0.008968117355803907
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1689346666666667
$ws.Range("H11").Value = 0.506804
$ws.Range("I11").Value = 0.06122536009783416
$ws.Range("J11").Value = 0.06122536009783416
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.799695333333333
$ws.Range("N11").Value = 17.399086
$ws.Range("O11").Value = 0.2811226515149324
$ws.Range("P11").Value = 0.2811226515149324
$ws.Range("Q11").Value = 0.9797695979048889
$ws.Range("R11").Value = 8.817926381144002
$ws.Range("S11").Value = 0.01721183557065968
$ws.Range("T11").Value = 0.01721183557065968
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1689346666666667
$ws.Range("H12").Value = 0.506804
$ws.Range("I12").Value = 0.06122536009783416
$ws.Range("J12").Value = 0.06122536009783416
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.259039333333333
$ws.Range("N12").Value = 18.777118
$ws.Range("O12").Value = 0.303387959572633
$ws.Range("P12").Value = 0.303387959572633
$ws.Range("Q12").Value = 1.057368723430222
$ws.Range("R12").Value = 9.516318510871999
$ws.Range("S12").Value = 0.0185750370741816
$ws.Range("T12").Value = 0.0185750370741816
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1689346666666667
$ws.Range("H13").Value = 0.506804
$ws.Range("I13").Value = 0.06122536009783416
$ws.Range("J13").Value = 0.06122536009783416
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.549851333333334
$ws.Range("N13").Value = 16.649554
$ws.Range("O13").Value = 0.2690122209305161
$ws.Range("P13").Value = 0.2690122209305161
$ws.Range("Q13").Value = 0.9375622850462223
$ws.Range("R13").Value = 8.438060565416002
$ws.Range("S13").Value = 0.01647037009718897
$ws.Range("T13").Value = 0.01647037009718897
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.4843330000000001
$ws.Range("H14").Value = 1.452999
$ws.Range("I14").Value = 0.1755321327313773
$ws.Range("J14").Value = 0.1755321327313773
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.021894333333333
$ws.Range("N14").Value = 9.065683
$ws.Range("O14").Value = 0.1464771679819186
$ws.Range("P14").Value = 0.1464771679819185
$ws.Range("Q14").Value = 1.463603148146334
$ws.Range("R14").Value = 13.172428333317
$ws.Range("S14").Value = 0.02571144969231838
$ws.Range("T14").Value = 0.02571144969231837
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.4843330000000001
$ws.Range("H15").Value = 1.452999
$ws.Range("I15").Value = 0.1755321327313773
$ws.Range("J15").Value = 0.1755321327313773
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.799695333333333
$ws.Range("N15").Value = 17.399086
$ws.Range("O15").Value = 0.2811226515149324
$ws.Range("P15").Value = 0.2811226515149324
$ws.Range("Q15").Value = 2.808983839879334
$ws.Range("R15").Value = 25.280854558914
$ws.Range("S15").Value = 0.04934605857951584
$ws.Range("T15").Value = 0.04934605857951584
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.4843330000000001
$ws.Range("H16").Value = 1.452999
$ws.Range("I16").Value = 0.1755321327313773
$ws.Range("J16").Value = 0.1755321327313773
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.259039333333333
$ws.Range("N16").Value = 18.777118
$ws.Range("O16").Value = 0.303387959572633
$ws.Range("P16").Value = 0.303387959572633
$ws.Range("Q16").Value = 3.031459297431334
$ws.Range("R16").Value = 27.283133676882
$ws.Range("S16").Value = 0.05325433558880514
$ws.Range("T16").Value = 0.05325433558880514
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.4843330000000001
$ws.Range("H17").Value = 1.452999
$ws.Range("I17").Value = 0.1755321327313773
$ws.Range("J17").Value = 0.1755321327313773
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.549851333333334
$ws.Range("N17").Value = 16.649554
$ws.Range("O17").Value = 0.2690122209305161
$ws.Range("P17").Value = 0.2690122209305161
$ws.Range("Q17").Value = 2.687976145827334
$ws.Range("R17").Value = 24.191785312446
$ws.Range("S17").Value = 0.04722028887073795
$ws.Range("T17").Value = 0.04722028887073794
